# Rerunning plots and creating latex table.
#
# Adds a "percent solved" statistic (row 547, cols J/K) to every
# benders_results_* sheet, and on the first sheet (benders_results_0)
# also adds cross-sheet "average ..." rollups (cols P/Q) next to the
# five headline statistic rows (percent solved, avg. gap,
# median iterations, avg. iteration time, avg. solve time).
# Finally, re-selects benders_results_0 as the active sheet/tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # benders_results_0
$ws2 = $wb.Worksheets.Item(2)   # benders_results_5
$ws3 = $wb.Worksheets.Item(3)   # benders_results_10
$ws4 = $wb.Worksheets.Item(4)   # benders_results_15

$pctFormula = "=100*COUNTIF(E20:E555,0)/ROWS(E20:E555)"

# --- every sheet gets the new "percent solved" row ---
foreach ($ws in @($ws1, $ws2, $ws3, $ws4)) {
    $ws.Range("J547").Value = "percent solved"
    $ws.Range("K547").Formula = $pctFormula
}

# --- benders_results_0 additionally gets the cross-sheet averages ---
$ws1.Range("P547").Value = "average percent solved"
$ws1.Range("Q547").Formula = "=AVERAGE(benders_results_0:benders_results_15!K547:K547)"

$ws1.Range("P549").Value = "average avg gap"
$ws1.Range("Q549").Formula = "=AVERAGE(benders_results_0:benders_results_15!K549:K549)"

$ws1.Range("P552").Value = "average median iterations"
$ws1.Range("Q552").Formula = "=AVERAGE(benders_results_0:benders_results_15!K552:K552)"

$ws1.Range("P553").Value = "average avg iteration time"
$ws1.Range("Q553").Formula = "=AVERAGE(benders_results_0:benders_results_15!K553:K553)"

$ws1.Range("P554").Value = "average avg solve time"
$ws1.Range("Q554").Formula = "=AVERAGE(benders_results_0:benders_results_15!K554:K554)"

# --- make benders_results_0 the selected / active tab again ---
$ws1.Activate()
$ws1.Range("Q552").Select()

$ws2.Range("P547").Select()
$ws3.Range("J547").Select()
$ws4.Range("J550").Select()

$ws1.Activate()
